# "add Q Read to Excell" -- duplicate the header/sample rows from Sheet1 into
# a new Sheet2, tweak a couple of cells, and leave Sheet2 as the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: update the selection / tab state -------------------------
[void]$ws1.Range("A1:D6").Select()

# --- Sheet2: new sheet, placed right after Sheet1 ----------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Bring over the header row plus the first few sample rows (values + the
# shared style) from Sheet1 in one shot so styles/shared-strings line up.
[void]$ws1.Range("A1:D6").Copy($ws2.Range("A1"))

# Row 2 becomes fully blank (was just "A" in A2 on Sheet1).
$ws2.Range("A2").ClearContents()

# Row 4: clear out everything except B4, which gets a brand new value.
$ws2.Range("A4").ClearContents()
$ws2.Range("C4").ClearContents()
$ws2.Range("D4").ClearContents()
$ws2.Range("B4").Value = "BBB"

# Extra trailing row with a single formatted-but-empty cell at B7 (copy
# the style from an already-blank, same-styled cell on Sheet1).
[void]$ws1.Range("B2").Copy($ws2.Range("B7"))

# Sizing to match the new sheet's look.
$ws2.Rows("1:6").RowHeight = 45
$ws2.Rows("7:7").RowHeight = 19
$ws2.Columns("A:E").ColumnWidth = 23.67

# Selection / active-tab state for the new sheet.
[void]$ws2.Range("A13").Select()

# Page setup tweak picked up on Sheet1 as part of the edit.
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1
